{"js": "// Remove the \"{{ }}\" mustache wrapper around each placeholder token in the\n// invoice template, e.g. \"{{INVOICE_NUMBER}}\" -> \"INVOICE_NUMBER\".\nconst placeholders = [\n  \"INVOICE_NUMBER\",\n  \"INVOICE_DATE\",\n  \"DUE_DATE\",\n  \"CLIENT_NAME\",\n  \"CLIENT_ADDRESS\",\n  \"CLIENT_EMAIL\",\n  \"ITEM_DESCRIPTION\",\n  \"QUANTITY\",\n  \"UNIT_PRICE\",\n  \"TOTAL_AMOUNT\",\n  \"PAYMENT_TERMS\",\n  \"NOTES\"\n];\n\nconst body = context.document.body;\n\nfor (const name of placeholders) {\n  const results = body.search(\"{{\" + name + \"}}\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(name, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Remove the \"{{ }}\" mustache wrapper around each placeholder token in the\n# invoice template, e.g. \"{{INVOICE_NUMBER}}\" -> \"INVOICE_NUMBER\".\n$d = $word.ActiveDocument\n\n$placeholders = @(\n    \"INVOICE_NUMBER\",\n    \"INVOICE_DATE\",\n    \"DUE_DATE\",\n    \"CLIENT_NAME\",\n    \"CLIENT_ADDRESS\",\n    \"CLIENT_EMAIL\",\n    \"ITEM_DESCRIPTION\",\n    \"QUANTITY\",\n    \"UNIT_PRICE\",\n    \"TOTAL_AMOUNT\",\n    \"PAYMENT_TERMS\",\n    \"NOTES\"\n)\n\nforeach ($name in $placeholders) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = \"{{\" + $name + \"}}\"\n    $find.MatchWildcards = $false\n    $find.Replacement.Text = $name\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
